$wb = $excel.ActiveWorkbook

# Rename sheet tabs: replace the first hyphen after "TRE" with an underscore
$wb.Worksheets.Item("Include from TRE-R67-TypeStru").Name = "Include from TRE_R67-TypeStru"
$wb.Worksheets.Item("Include from TRE-R04-TypeSavo").Name = "Include from TRE_R04-TypeSavo"
$wb.Worksheets.Item("Include from TRE-R288-TypePro").Name = "Include from TRE_R288-TypePro"

# Update the Date value on the Metadata sheet
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value = "2024-04-08T14:06:04+00:00"
